$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin names, URLs) - safe to assign directly ---
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"

# --- Numeric-looking / percentage text cells (Price, Volume columns) ---
# Force text storage via NumberFormat "@" then restore default style,
# so Excel does not silently convert these into real numbers.
$deCells = @{
    "D2" = "72.335.75"
    "E2" = "  +0.19%  "
    "D3" = "2.632.84"
    "E3" = "  -1.26%  "
    "E4" = "  -0.01%  "
    "D5" = "585.33"
    "E5" = "  -2.19%  "
    "D6" = "175.07"
    "E6" = "  -0.45%  "
    "E7" = "  -0.04%  "
    "E8" = "  -0.50%  "
    "D9" = "0.172"
    "E9" = "  +1.94%  "
    "D10" = "2.632.19"
    "E10" = "  -1.20%  "
    "E11" = "  +1.46%  "
    "D12" = "0.358"
    "E12" = "  +1.61%  "
    "E13" = "  -1.69%  "
    "D14" = "3.115.07"
    "E14" = "  -1.30%  "
    "E15" = "  +0.57%  "
    "D16" = "72.194.73"
    "E16" = "  +0.14%  "
    "D17" = "25.71"
    "E17" = "  -1.88%  "
    "D18" = "2.631.01"
    "E18" = "  -1.41%  "
    "E19" = "  +0.29%  "
    "D20" = "7.85"
    "E20" = "  -1.44%  "
    "D21" = "374.90"
    "E21" = "  +1.37%  "
    "E22" = "  -1.24%  "
    "E23" = "  -0.02%  "
    "D24" = "71.49"
    "E24" = "  -0.18%  "
    "E25" = "  -0.08%  "
    "D26" = "4.22"
    "E26" = "  -2.28%  "
    "D27" = "9.47"
    "E27" = "  -3.23%  "
    "D28" = "2.769.30"
    "E28" = "  -1.24%  "
    "D29" = "1.00"
    "E29" = "  +0.08%  "
    "D30" = "0.0₃0946"
    "E30" = "  +1.18%  "
    "D31" = "7.94"
    "E31" = "  -1.16%  "
    "D32" = "490.81"
    "E32" = "  -3.48%  "
    "E34" = "  -0.82%  "
    "D35" = "0.999"
    "E35" = "  +0.04%  "
    "D36" = "161.43"
    "E36" = "  -2.03%  "
    "E37" = "  +9.44%  "
    "D38" = "19.14"
    "E38" = "  -1.89%  "
    "E39" = "  -1.07%  "
    "E40" = "  -0.83%  "
    "E41" = "  -0.14%  "
    "D42" = "1.73"
    "E42" = "  -4.69%  "
    "D43" = "2.56"
    "E43" = "  +0.15%  "
    "E44" = "  -2.29%  "
    "E45" = "  -1.93%  "
    "D46" = "39.04"
    "E46" = "  -0.47%  "
    "D47" = "150.16"
    "E47" = "  -1.63%  "
    "E48" = "  -2.46%  "
    "E49" = "  -0.98%  "
    "E50" = "  -2.65%  "
    "D51" = "0.605"
    "E51" = "  +0.94%  "
}
foreach ($ref in $deCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $deCells[$ref]
    $cell.Style = "Normal"
}
